# Update the date line in the first paragraph.
$d = $word.ActiveDocument

$mult = [char]0x00D7   # the "x" (multiplication sign, U+00D7) used in the table cells

$d.Paragraphs.Item(1).Range.Find.Execute(
    "2025-05-28 Wednesday", $false, $false, $false, $false, $false,
    $true, 1, $false, "2025-05-29 Thursday", 2) | Out-Null

# Update the multiplication exercises held in the table cells. Each data row
# of the table (1, 5, 10, 15, 20) holds five exercises, one per column, and
# is addressed directly by row/column so that duplicate problem text (e.g.
# "634x3=" appearing twice) is handled unambiguously.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "238$($mult)8=" },
    @{ Row = 1;  Col = 2; Text = "234$($mult)4=" },
    @{ Row = 1;  Col = 3; Text = "979$($mult)9=" },
    @{ Row = 1;  Col = 4; Text = "456$($mult)5=" },
    @{ Row = 1;  Col = 5; Text = "182$($mult)3=" },

    @{ Row = 5;  Col = 1; Text = "129$($mult)6=" },
    @{ Row = 5;  Col = 2; Text = "343$($mult)6=" },
    @{ Row = 5;  Col = 3; Text = "356$($mult)2=" },
    @{ Row = 5;  Col = 4; Text = "280$($mult)3=" },
    @{ Row = 5;  Col = 5; Text = "935$($mult)5=" },

    @{ Row = 10; Col = 1; Text = "174$($mult)3=" },
    @{ Row = 10; Col = 2; Text = "345$($mult)6=" },
    @{ Row = 10; Col = 3; Text = "144$($mult)8=" },
    @{ Row = 10; Col = 4; Text = "167$($mult)2=" },
    @{ Row = 10; Col = 5; Text = "234$($mult)7=" },

    @{ Row = 15; Col = 1; Text = "120$($mult)2=" },
    @{ Row = 15; Col = 2; Text = "233$($mult)8=" },
    @{ Row = 15; Col = 3; Text = "737$($mult)3=" },
    @{ Row = 15; Col = 4; Text = "292$($mult)8=" },
    @{ Row = 15; Col = 5; Text = "803$($mult)9=" },

    @{ Row = 20; Col = 1; Text = "262$($mult)2=" },
    @{ Row = 20; Col = 2; Text = "990$($mult)4=" },
    @{ Row = 20; Col = 3; Text = "523$($mult)5=" },
    @{ Row = 20; Col = 4; Text = "214$($mult)3=" },
    @{ Row = 20; Col = 5; Text = "785$($mult)4=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
